$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H47").Value = 4500
$ws.Range("I47").Value = 4500
$ws.Range("K47").Value = 4500
$ws.Range("M47").Value = -3528
$ws.Range("H54").Value = 25042
$ws.Range("I54").Value = 20000
$ws.Range("J54").Value = 30084
$ws.Range("K54").Value = 20000
$ws.Range("L54").Value = 30084
$ws.Range("M54").Value = -19514
$ws.Range("N54").Value = -31056
$ws.Range("H58").Value = 111.125
$ws.Range("J58").Value = 350
$ws.Range("L58").Value = 1050
$ws.Range("N58").Value = -1350
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("H98").Value = 937.0909
$ws.Range("I98").Value = 937.0909
$ws.Range("K98").Value = 937.0909
$ws.Range("M98").Value = 560.9091
$ws.Range("H106").Value = 4857.143
$ws.Range("I106").Value = 5000
$ws.Range("J106").Value = 4000
$ws.Range("K106").Value = 5000
$ws.Range("L106").Value = 4000
$ws.Range("M106").Value = -4369
$ws.Range("N106").Value = -5262
$ws.Range("H122").Value = 937.0909
$ws.Range("I122").Value = 937.0909
$ws.Range("K122").Value = 2811.2727
$ws.Range("M122").Value = -361.2727
$ws.Range("H137").Value = 3174.577
$ws.Range("I137").Value = 2044.6666
$ws.Range("J137").Value = 4143.0713
$ws.Range("K137").Value = 6133.9998
$ws.Range("L137").Value = 12429.2139
$ws.Range("M137").Value = -3583.9998
$ws.Range("N137").Value = -17529.2139
$ws.Range("N76").ClearContents()
$ws.Range("N79").ClearContents()

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2462.4211
$ws.Range("I32").Value = 2462.4211
$ws.Range("K32").Value = 2462.4211
$ws.Range("M32").Value = -2175.4211
$ws.Range("H44").Value = 11614.35
$ws.Range("J44").Value = 11614.35
$ws.Range("L44").Value = 11614.35
$ws.Range("N44").Value = -12590.35
$ws.Range("H45").Value = 2751.12
$ws.Range("I45").Value = 2375.1428
$ws.Range("J45").Value = 4725
$ws.Range("K45").Value = 2375.1428
$ws.Range("L45").Value = 4725
$ws.Range("M45").Value = -1998.1428
$ws.Range("N45").Value = -5479
$ws.Range("H55").Value = 41833
$ws.Range("J55").Value = 58749.5
$ws.Range("L55").Value = 58749.5
$ws.Range("N55").Value = -59379.5
$ws.Range("H97").Value = 540.3684
$ws.Range("I97").Value = 540.3684
$ws.Range("K97").Value = 540.3684
$ws.Range("M97").Value = -44.36839999999995
$ws.Range("H122").Value = 1426.5385
$ws.Range("I122").Value = 1426.5385
$ws.Range("K122").Value = 4279.6155
$ws.Range("M122").Value = -1829.6155
$ws.Range("H124").Value = 76666.336
$ws.Range("J124").Value = 76666.336
$ws.Range("L124").Value = 76666.336
$ws.Range("N124").Value = -86486.336
$ws.Range("H125").Value = 91125
$ws.Range("J125").Value = 91125
$ws.Range("L125").Value = 91125
$ws.Range("N125").Value = -100965
$ws.Range("H132").Value = 2711.4285
$ws.Range("I132").Value = 1217.7778
$ws.Range("J132").Value = 5400
$ws.Range("K132").Value = 3653.3334
$ws.Range("L132").Value = 16200
$ws.Range("M132").Value = -1123.3334
$ws.Range("N132").Value = -21260
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N138").ClearContents()
$ws.Range("N139").ClearContents()

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 8560.666999999999
$ws.Range("I20").Value = 8166
$ws.Range("K20").Value = 8166
$ws.Range("M20").Value = -7919
$ws.Range("H25").Value = 4007
$ws.Range("I25").Value = 4007
$ws.Range("K25").Value = 4007
$ws.Range("M25").Value = -3772
$ws.Range("H134").Value = 7127
$ws.Range("I134").Value = 6552.4
$ws.Range("K134").Value = 19657.2
$ws.Range("M134").Value = -17122.2

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6182.9062
$ws.Range("I31").Value = 2729.5625
$ws.Range("J31").Value = 9636.25
$ws.Range("K31").Value = 2729.5625
$ws.Range("L31").Value = 9636.25
$ws.Range("M31").Value = -2434.5625
$ws.Range("N31").Value = -10226.25
$ws.Range("H34").Value = 6182.9062
$ws.Range("I34").Value = 2729.5625
$ws.Range("J34").Value = 9636.25
$ws.Range("K34").Value = 2729.5625
$ws.Range("L34").Value = 9636.25
$ws.Range("M34").Value = -2527.5625
$ws.Range("N34").Value = -10040.25
$ws.Range("H122").Value = 1834.9584
$ws.Range("I122").Value = 1636.5
$ws.Range("K122").Value = 4909.5
$ws.Range("M122").Value = -2459.5

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 78.454544
$ws.Range("I7").Value = 76.09999999999999
$ws.Range("J7").Value = 102
$ws.Range("K7").Value = 228.3
$ws.Range("L7").Value = 306
$ws.Range("M7").Value = -116.3
$ws.Range("N7").Value = -530
$ws.Range("H8").Value = 350.75
$ws.Range("I8").Value = 350.75
$ws.Range("K8").Value = 1052.25
$ws.Range("M8").Value = -913.25
$ws.Range("H98").Value = 234.8
$ws.Range("I98").Value = 246
$ws.Range("J98").Value = 190
$ws.Range("K98").Value = 738
$ws.Range("L98").Value = 570
$ws.Range("M98").Value = 760
$ws.Range("N98").Value = -3566
$ws.Range("H107").Value = 1400
$ws.Range("J107").Value = 1400
$ws.Range("L107").Value = 4200
$ws.Range("N107").Value = -8040
$ws.Range("H108").Value = 1504.6
$ws.Range("I108").Value = 924.6667
$ws.Range("K108").Value = 2774.0001
$ws.Range("M108").Value = 105.9998999999998
$ws.Range("H109").Value = 168835.67
$ws.Range("I109").Value = 251856
$ws.Range("J109").Value = 2795
$ws.Range("K109").Value = 755568
$ws.Range("L109").Value = 8385
$ws.Range("M109").Value = -754528
$ws.Range("N109").Value = -10465
$ws.Range("H131").Value = 993.5
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("H134").Value = 9411.143
$ws.Range("I134").Value = 969.75
$ws.Range("K134").Value = 2909.25
$ws.Range("M134").Value = 2160.75
$ws.Range("N131").ClearContents()

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 26000
$ws.Range("I46").Value = 20000
$ws.Range("J46").Value = 32000
$ws.Range("K46").Value = 20000
$ws.Range("L46").Value = 32000
$ws.Range("M46").Value = -19844
$ws.Range("N46").Value = -32312
$ws.Range("H80").Value = 3474
$ws.Range("I80").Value = 3299.3333
$ws.Range("J80").Value = 3998
$ws.Range("K80").Value = 3299.3333
$ws.Range("L80").Value = 3998
$ws.Range("M80").Value = -2301.3333
$ws.Range("N80").Value = -5994
$ws.Range("H83").Value = 3474
$ws.Range("I83").Value = 3299.3333
$ws.Range("J83").Value = 3998
$ws.Range("K83").Value = 16496.6665
$ws.Range("L83").Value = 19990
$ws.Range("M83").Value = -11504.6665
$ws.Range("N83").Value = -29974
$ws.Range("H122").Value = 298281.8
$ws.Range("I122").Value = 337203.6
$ws.Range("K122").Value = 1011610.8
$ws.Range("M122").Value = -1009160.8

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3383.476
$ws.Range("I61").Value = 2297.2354
$ws.Range("K61").Value = 2297.2354
$ws.Range("M61").Value = -2095.2354
$ws.Range("H100").Value = 6147.421
$ws.Range("I100").Value = 2731
$ws.Range("K100").Value = 2731
$ws.Range("M100").Value = -2190
$ws.Range("H113").Value = 3383.476
$ws.Range("I113").Value = 2297.2354
$ws.Range("K113").Value = 2297.2354
$ws.Range("M113").Value = -127.2354
$ws.Range("H122").Value = 4418.75
$ws.Range("I122").Value = 3905.3333
$ws.Range("K122").Value = 11715.9999
$ws.Range("M122").Value = -9265.999899999999
$ws.Range("H132").Value = 5998
$ws.Range("I132").Value = 5196.6
$ws.Range("K132").Value = 15589.8
$ws.Range("M132").Value = -13059.8

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 4227.278
$ws.Range("I126").Value = 1871.5454
$ws.Range("J126").Value = 7929.143
$ws.Range("K126").Value = 5614.6362
$ws.Range("L126").Value = 23787.429
$ws.Range("M126").Value = -3144.6362
$ws.Range("N126").Value = -28727.429
$ws.Range("H140").Value = 0
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("M140").ClearContents()
$ws.Range("N140").ClearContents()
